$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.243.62'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '2.296.95'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '317.86'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').Value = '102.71'
$ws.Range('E6').Value = '  -3.60%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.606'
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('D10').Value = '39.64'
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('D12').Value = '8.41'
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').Value = '0.957'
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').Value = '15.27'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').Value = '2.644.26'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').Value = '2.296.65'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').Value = '42.259.96'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('E19').Value = '  -1.48%  '
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('D21').Value = '12.56'
$ws.Range('E21').Value = '  +34.01%  '
$ws.Range('D22').Value = '73.42'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('E23').Value = '  +2.40%  '
$ws.Range('D24').Value = '275.85'
$ws.Range('E24').Value = '  +7.42%  '
$ws.Range('D25').Value = '2.25'
$ws.Range('E25').Value = '  -2.78%  '
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('D27').Value = '10.82'
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').Value = '2.38'
$ws.Range('E28').Value = '  +6.84%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').Value = '37.62'
$ws.Range('E30').Value = '  +5.51%  '
$ws.Range('D31').Value = '165.66'
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('D32').Value = '6.08'
$ws.Range('E32').Value = '  +5.08%  '
$ws.Range('D33').Value = '0.0875'
$ws.Range('E33').Value = '  -2.02%  '
$ws.Range('E34').Value = '  +3.30%  '
$ws.Range('D35').Value = '2.66'
$ws.Range('E35').Value = '  -8.68%  '
$ws.Range('D36').Value = '0.118'
$ws.Range('E36').Value = '  -0.81%  '
$ws.Range('D37').Value = '4.57'
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('E38').Value = '  +1.76%  '
$ws.Range('D39').Value = '3.70'
$ws.Range('E39').Value = '  +2.04%  '
$ws.Range('E40').Value = '  -4.16%  '
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').Value = '69.71'
$ws.Range('E42').Value = '  -2.80%  '
$ws.Range('B43').Value = 'BitcoinSV'
$ws.Range('C43').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D43').Value = '95.71'
$ws.Range('E43').Value = '  -1.15%  '
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = '11.99'
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('D47').Value = '112.53'
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('D48').Value = '79.17'
$ws.Range('E48').Value = '  +5.17%  '
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('D50').Value = '5.26'
$ws.Range('E50').Value = '  -1.01%  '
$ws.Range('D51').Value = '1.596.56'
$ws.Range('E51').Value = '  +3.09%  '
